$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 13: Emilio Rugerio, second SmartScore submission (20251128_172733) ---

# Pre-format the SmartScore cells for this row as Text so the
# trailing-zero score strings (e.g. "0.580") are preserved verbatim
# instead of being coerced to numbers.
$smartScoreCells = @("I13","L13","O13","R13","U13","X13","AA13","AD13","AG13")
foreach ($cellAddr in $smartScoreCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

$ws.Range("A13").Value = @'
Emilio Rugerio_20251128_172733
'@

$ws.Range("B13").Value = ""

$ws.Range("C13").Value = @'
Emilio Rugerio
'@

$ws.Range("D13").Value = 21

$ws.Range("E13").Value = @'
Male
'@

$ws.Range("F13").Value = @'
2025-11-28 17:27:33
'@

$ws.Range("G13").Value = @'
{
  "portion": 0.6,
  "diet": 0.2857142857142857,
  "salt": 0.2,
  "fat": 0.2,
  "natural": 0.2,
  "convenience": 0.2,
  "price": 0.4
}
'@

$ws.Range("H13").Value = @'
Nongshim Neoguri Spicy Seafood
'@

$ws.Range("I13").Value = @'
0.642
'@

$ws.Range("J13").Value = @'
Sabor a marisco, umami, picante equilibrado, buena textura, algo salado
'@

$ws.Range("K13").Value = @'
Nissin Chow Mein Teriyaki Beef
'@

$ws.Range("L13").Value = @'
0.615
'@

$ws.Range("M13").Value = @'
Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa
'@

$ws.Range("N13").Value = @'
Nongshim Shin Ramyun
'@

$ws.Range("O13").Value = @'
0.562
'@

$ws.Range("P13").Value = @'
Sabor intenso, picante, umami, fideos gruesos, muy alto en sodio
'@

$ws.Range("Q13").Value = @'
Amy’s Macaroni & Cheese (frozen)
'@

$ws.Range("R13").Value = @'
0.601
'@

$ws.Range("S13").Value = @'
Queso real, textura casera, sin conservadores, alto en grasa, algo caro
'@

$ws.Range("T13").Value = @'
Kraft Macaroni & Cheese Dinner
'@

$ws.Range("U13").Value = @'
0.580
'@

$ws.Range("V13").Value = @'
Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato
'@

$ws.Range("W13").Value = @'
Annie’s Shells & White Cheddar
'@

$ws.Range("X13").Value = @'
0.517
'@

$ws.Range("Y13").Value = @'
Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños
'@

$ws.Range("Z13").Value = @'
Wild Planet Wild Tuna Pasta Salad
'@

$ws.Range("AA13").Value = @'
0.621
'@

$ws.Range("AB13").Value = @'
Sabor fresco, buena proteína, saludable, porción algo pequeña
'@

$ws.Range("AC13").Value = @'
StarKist Chicken Creations (Chicken Salad)
'@

$ws.Range("AD13").Value = @'
0.507
'@

$ws.Range("AE13").Value = @'
Portátil, saludable, fácil, buena textura, sabor suave
'@

$ws.Range("AF13").Value = @'
Jack Link’s Beef Jerky Original
'@

$ws.Range("AG13").Value = @'
0.487
'@

$ws.Range("AH13").Value = @'
Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña
'@
